$d = $word.ActiveDocument

$d.Content.Find.Execute("Tela09_TelaEntrada", $true, $false, $false, $false, $false, $true, 1, $false, "Tela05_TelaEntrada", 2)
$d.Content.Find.Execute("Tela10_perfilOpcoes", $true, $false, $false, $false, $false, $true, 1, $false, "Tela06_perfilOpcoes", 2)
$d.Content.Find.Execute("Tela13_HistoricoDeEventos", $true, $false, $false, $false, $false, $true, 1, $false, "Tela09_HistoricoDeEventos", 2)
$d.Content.Find.Execute("Tela20_StatusDaInscrição", $true, $false, $false, $false, $false, $true, 1, $false, "Tela14_StatusDaInscrição", 2)
